$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 09:52"

# Row 18 -> Rusia updated stats
$ws.Cells.Item(18,2).Value = 24490
$ws.Cells.Item(18,3).Value = 3388
$ws.Cells.Item(18,4).Value = 1986
$ws.Cells.Item(18,5).Value = 22306
$ws.Cells.Item(18,6).Value = 8
$ws.Cells.Item(18,7).Value = 28
$ws.Cells.Item(18,8).Value = 198

# Row 69 -> Kazajistan updated stats
$ws.Cells.Item(69,2).Value = 1275
$ws.Cells.Item(69,3).Value = 43
$ws.Cells.Item(69,4).Value = 203
$ws.Cells.Item(69,5).Value = 1057
$ws.Cells.Item(69,6).Value = 20
$ws.Cells.Item(69,7).Value = 1
$ws.Cells.Item(69,8).Value = 15

# Row 73 -> Armenia
$ws.Cells.Item(73,1).Value = "Armenia"
$ws.Cells.Item(73,2).Value = 1111
$ws.Cells.Item(73,3).Value = 44
$ws.Cells.Item(73,4).Value = 297
$ws.Cells.Item(73,5).Value = 797
$ws.Cells.Item(73,6).Value = 30
$ws.Cells.Item(73,7).Value = 1
$ws.Cells.Item(73,8).Value = 17

# Row 74 -> Lituania
$ws.Cells.Item(74,1).Value = "Lituania"
$ws.Cells.Item(74,2).Value = 1091
$ws.Cells.Item(74,3).Value = 21
$ws.Cells.Item(74,4).Value = 138
$ws.Cells.Item(74,5).Value = 924
$ws.Cells.Item(74,6).Value = 14
$ws.Cells.Item(74,7).Value = 0
$ws.Cells.Item(74,8).Value = 29

# Row 75 -> Bosnia y Herzegovina
$ws.Cells.Item(75,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(75,2).Value = 1083
$ws.Cells.Item(75,3).Value = 0
$ws.Cells.Item(75,4).Value = 236
$ws.Cells.Item(75,5).Value = 807
$ws.Cells.Item(75,6).Value = 4
$ws.Cells.Item(75,7).Value = 0
$ws.Cells.Item(75,8).Value = 40
